$wb = $excel.ActiveWorkbook

# Worksheets involved in the edit
$wsBuscaLupa = $wb.Worksheets.Item("BuscaLupa")

# Swap the data values in row 1 (A1 <-> B1) and row 2 (A2 <-> B2)
$a1 = $wsBuscaLupa.Range("A1").Value2
$b1 = $wsBuscaLupa.Range("B1").Value2
$wsBuscaLupa.Range("A1").Value2 = $b1
$wsBuscaLupa.Range("B1").Value2 = $a1

$a2 = $wsBuscaLupa.Range("A2").Value2
$b2 = $wsBuscaLupa.Range("B2").Value2
$wsBuscaLupa.Range("A2").Value2 = $b2
$wsBuscaLupa.Range("B2").Value2 = $a2

# Make BuscaLupa the active sheet/tab and move the selection to B2
$wsBuscaLupa.Activate()
$wsBuscaLupa.Range("B2").Select()

$wb.Save()
